$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.197282791137695
$ws.Range("B1").Value = 2.60010552406311
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.184739589691162
$ws.Range("E1").Value = 1.176159024238586
